$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3044.2727
$ws.Range("I53").Value = 355.875
$ws.Range("J53").Value = 10213.333
$ws.Range("K53").Value = 355.875
$ws.Range("L53").Value = 10213.333
$ws.Range("M53").Value = 281.125
$ws.Range("N53").Value = -11487.333

$ws.Range("H129").Value = 209446.48
$ws.Range("J129").Value = 245154.95
$ws.Range("L129").Value = 735464.8500000001
$ws.Range("N129").Value = -745464.8500000001

$ws.Range("H132").Value = 2493.8667
$ws.Range("I132").Value = 2519.0527
$ws.Range("J132").Value = 2357.1428
$ws.Range("K132").Value = 7557.158100000001
$ws.Range("L132").Value = 7071.428400000001
$ws.Range("M132").Value = -5027.158100000001
$ws.Range("N132").Value = -12131.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H122").Value = 1838.3478
$ws.Range("I122").Value = 1775.15
$ws.Range("J122").Value = 2259.6667
$ws.Range("K122").Value = 5325.450000000001
$ws.Range("L122").Value = 6779.000100000001
$ws.Range("M122").Value = -2875.450000000001
$ws.Range("N122").Value = -11679.0001

$ws.Range("H124").Value = 11485
$ws.Range("J124").Value = 11485
$ws.Range("L124").Value = 11485
$ws.Range("N124").Value = -21305

$ws.Range("H125").Value = 28315.125
$ws.Range("J125").Value = 28315.125
$ws.Range("L125").Value = 28315.125
$ws.Range("N125").Value = -38155.125

$ws.Range("H132").Value = 9488.538
$ws.Range("I132").Value = 1747.4043
$ws.Range("J132").Value = 29701.5
$ws.Range("K132").Value = 5242.2129
$ws.Range("L132").Value = 89104.5
$ws.Range("M132").Value = -2712.2129
$ws.Range("N132").Value = -94164.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1639.7778
$ws.Range("I20").Value = 1582.4
$ws.Range("K20").Value = 1582.4
$ws.Range("M20").Value = -1335.4

$ws.Range("H94").Value = 904.25
$ws.Range("I94").Value = 722.3333
$ws.Range("J94").Value = 1450
$ws.Range("K94").Value = 722.3333
$ws.Range("L94").Value = 1450
$ws.Range("M94").Value = -271.3333
$ws.Range("N94").Value = -2352

$ws.Range("H99").Value = 1594.25
$ws.Range("I99").Value = 1705
$ws.Range("J99").Value = 1483.5
$ws.Range("K99").Value = 1705
$ws.Range("L99").Value = 1483.5
$ws.Range("M99").Value = -207
$ws.Range("N99").Value = -4479.5

$ws.Range("H107").Value = 1389.8572
$ws.Range("I107").Value = 1305.0435
$ws.Range("K107").Value = 1305.0435
$ws.Range("M107").Value = 614.9565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3394.1924
$ws.Range("J99").Value = 5800
$ws.Range("L99").Value = 5800
$ws.Range("N99").Value = -8796

$ws.Range("H126").Value = 3394.1924
$ws.Range("J126").Value = 5800
$ws.Range("L126").Value = 17400
$ws.Range("N126").Value = -22340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 470
$ws.Range("J7").Value = 575
$ws.Range("L7").Value = 1725
$ws.Range("N7").Value = -1949

$ws.Range("H92").Value = 374.83334
$ws.Range("J92").Value = 429.8
$ws.Range("L92").Value = 1289.4
$ws.Range("N92").Value = -3785.4

$ws.Range("H122").Value = 919
$ws.Range("J122").Value = 1022.0571
$ws.Range("L122").Value = 9198.5139
$ws.Range("N122").Value = -14098.5139

$ws.Range("H131").Value = 711.4
$ws.Range("J131").Value = 728.26044
$ws.Range("L131").Value = 2184.78132
$ws.Range("N131").Value = -12264.78132

$ws.Range("H136").Value = 2721
$ws.Range("I136").Value = 1202.2222
$ws.Range("J136").Value = 4999.1665
$ws.Range("K136").Value = 3606.6666
$ws.Range("L136").Value = 14997.4995
$ws.Range("M136").Value = 1493.3334
$ws.Range("N136").Value = -25197.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26534
$ws.Range("J57").Value = 26534
$ws.Range("L57").Value = 26534
$ws.Range("N57").Value = -28174

$ws.Range("H70").Value = 2984502.5
$ws.Range("I70").Value = 3807.2856
$ws.Range("K70").Value = 3807.2856
$ws.Range("M70").Value = -3537.2856

$ws.Range("H73").Value = 2984502.5
$ws.Range("I73").Value = 3807.2856
$ws.Range("K73").Value = 3807.2856
$ws.Range("M73").Value = -2871.2856

$ws.Range("H122").Value = 4159.067
$ws.Range("I122").Value = 3262.3635
$ws.Range("K122").Value = 9787.0905
$ws.Range("M122").Value = -7337.0905

$ws.Range("H126").Value = 3937.5217
$ws.Range("I126").Value = 4688.1763
$ws.Range("J126").Value = 3497.4827
$ws.Range("K126").Value = 14064.5289
$ws.Range("L126").Value = 10492.4481
$ws.Range("M126").Value = -11594.5289
$ws.Range("N126").Value = -15432.4481

$ws.Range("H132").Value = 33757.332
$ws.Range("I132").Value = 7393.1665
$ws.Range("J132").Value = 86485.664
$ws.Range("K132").Value = 22179.4995
$ws.Range("L132").Value = 259456.992
$ws.Range("M132").Value = -19649.4995
$ws.Range("N132").Value = -264516.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4466.7144
$ws.Range("I61").Value = 2294.5
$ws.Range("J61").Value = 17500
$ws.Range("K61").Value = 2294.5
$ws.Range("L61").Value = 17500
$ws.Range("M61").Value = -2092.5
$ws.Range("N61").Value = -17904

$ws.Range("H113").Value = 4466.7144
$ws.Range("I113").Value = 2294.5
$ws.Range("J113").Value = 17500
$ws.Range("K113").Value = 2294.5
$ws.Range("L113").Value = 17500
$ws.Range("M113").Value = -124.5
$ws.Range("N113").Value = -21840

$ws.Range("H122").Value = 1311126.6
$ws.Range("I122").Value = 1637460.4
$ws.Range("K122").Value = 4912381.199999999
$ws.Range("M122").Value = -4909931.199999999

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 2355.7144
$ws.Range("I132").Value = 1678.4
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 5035.200000000001
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -2505.200000000001
$ws.Range("N132").Value = -29060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
